$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 293, shifting existing rows 293:307 down to 294:308
$ws.Rows.Item(293).Insert()

# Populate the newly inserted row 293 with the new weekly data record
$ws.Cells.Item(293, 1).Value = 5
$ws.Cells.Item(293, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(293, 3).Value = "Maule"
$ws.Cells.Item(293, 4).Value = 44753
$ws.Cells.Item(293, 4).NumberFormat = $ws.Cells.Item(294, 4).NumberFormat
$ws.Cells.Item(293, 5).Value = 7
$ws.Cells.Item(293, 6).Value = 100112006
$ws.Cells.Item(293, 7).Value = "Repollo"
$ws.Cells.Item(293, 8).Value = "Crespo record"
$ws.Cells.Item(293, 9).Value = "Primera"
$ws.Cells.Item(293, 10).Value = 3000
$ws.Cells.Item(293, 11).Value = 1100
$ws.Cells.Item(293, 12).Value = 1100
$ws.Cells.Item(293, 13).Value = 1100
$ws.Cells.Item(293, 14).Value = "`$/unidad"
$ws.Cells.Item(293, 15).Value = "Región del Maule"
$ws.Cells.Item(293, 16).Value = 1100
$ws.Cells.Item(293, 17).Value = 1
$ws.Cells.Item(293, 18).Value = "Hortaliza"
